# Apply updated "想去人数" (want-to-go count) values to 展览 and 全部类型 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 201
$ws.Range("F4").Value = 5248
$ws.Range("F6").Value = 52
$ws.Range("F8").Value = 590
$ws.Range("F9").Value = 548
$ws.Range("F10").Value = 1047
$ws.Range("F12").Value = 1450
$ws.Range("F13").Value = 4192
$ws.Range("F14").Value = 432
$ws.Range("F15").Value = 173
$ws.Range("F16").Value = 156
$ws.Range("F17").Value = 94
$ws.Range("F18").Value = 3271
$ws.Range("F19").Value = 159
$ws.Range("F20").Value = 1072
$ws.Range("F21").Value = 100
$ws.Range("F23").Value = 189
$ws.Range("F24").Value = 112
$ws.Range("F25").Value = 33
$ws.Range("F26").Value = 135
$ws.Range("F27").Value = 69
$ws.Range("F28").Value = 294
$ws.Range("F30").Value = 54
$ws.Range("F32").Value = 18
$ws.Range("F33").Value = 18

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 201
$ws.Range("F5").Value = 5248
$ws.Range("F7").Value = 52
$ws.Range("F9").Value = 590
$ws.Range("F10").Value = 548
$ws.Range("F11").Value = 1047
$ws.Range("F13").Value = 1450
$ws.Range("F14").Value = 4193
$ws.Range("F15").Value = 432
$ws.Range("F16").Value = 173
$ws.Range("F17").Value = 156
$ws.Range("F18").Value = 94
$ws.Range("F19").Value = 3271
$ws.Range("F20").Value = 159
$ws.Range("F21").Value = 1072
$ws.Range("F22").Value = 100
$ws.Range("F24").Value = 189
$ws.Range("F25").Value = 112
$ws.Range("F26").Value = 33
$ws.Range("F27").Value = 135
$ws.Range("F28").Value = 69
$ws.Range("F29").Value = 294
$ws.Range("F31").Value = 54
$ws.Range("F33").Value = 18
$ws.Range("F34").Value = 18
